# Auto-generated edit script applying the diff to before.xlsx
# (regenerates the cell-level changes recorded in the commit diff)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 297
$ws.Range("F6").Value = 429
$ws.Range("F7").Value = 369
$ws.Range("F8").Value = 1944
$ws.Range("F10").Value = 27
$ws.Range("F11").Value = 23
$ws.Range("F12").Value = 1582
$ws.Range("F13").Value = 1582
$ws.Range("F14").Value = 1312
$ws.Range("F15").Value = 45
$ws.Range("F16").Value = 1381
$ws.Range("F18").Value = 10
$ws.Range("F19").Value = 1
$ws.Range("F20").Value = 425
$ws.Range("F23").Value = 134
$ws.Range("F24").Value = 6910
$ws.Range("F25").Value = 7472
$ws.Range("F26").Value = 30
$ws.Range("F29").Value = 44
$ws.Range("F30").Value = 214
$ws.Range("F31").Value = 234
$ws.Range("F32").Value = 9
$ws.Range("F33").Value = 5
$ws.Range("F36").Value = 1360
$ws.Range("F37").Value = 8
$ws.Range("F39").Value = 273
$ws.Range("F40").Value = 664
$ws.Range("F43").Value = 302
$ws.Range("F45").Value = 178
$ws.Range("F46").Value = 76
$ws.Range("F47").Value = 105
$ws.Range("F48").Value = 126
$ws.Range("F49").Value = 11

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 46
$ws.Range("F17").Value = 273

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 2535
$ws.Range("F5").Value = 104
$ws.Range("F6").Value = 2

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 104
$ws.Range("F9").Value = 297
$ws.Range("F11").Value = 429
$ws.Range("F12").Value = 369
$ws.Range("F13").Value = 1944
$ws.Range("F14").Value = 27
$ws.Range("F15").Value = 23
$ws.Range("F16").Value = 1582
$ws.Range("F17").Value = 1582
$ws.Range("C18").Value = "北京·不舍昼夜3.0-奇妙童话夜"
$ws.Range("D18").Value = "酒仙桥北路2号院798艺术区706后街1号 北京格瑞斯艺术酒店"
$ws.Range("E18").Value = "2024.09.16 20:00-09.17 02:00"
$ws.Range("F18").Value = 45
$ws.Range("G18").Value = 158
$ws.Range("H18").Value = "https://show.bilibili.com/platform/detail.html?id=91042"
$ws.Range("I18").Value = "//i0.hdslb.com/bfs/openplatform/202408/ZRwTjxgi1724204402969.jpeg"
$ws.Range("C19").Value = "北京·原神×星穹铁道only2.0同人展"
$ws.Range("D19").Value = "高碑店东路超级蜂巢 5G直播基地"
$ws.Range("E19").Value = "2024.09.16 10:00-09.16 17:00"
$ws.Range("F19").Value = 1381
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = "https://show.bilibili.com/platform/detail.html?id=88285"
$ws.Range("I19").Value = "//i1.hdslb.com/bfs/openplatform/202406/iWlE3Q9X1719554169582.jpeg"
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "2024-09-21"
$ws.Range("B20").NumberFormat = "General"
$ws.Range("C20").Value = "北京·核聚变游戏嘉年华2024"
$ws.Range("D20").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws.Range("E20").Value = "2024.09.21 09:00-09.22 17:00"
$ws.Range("F20").Value = 425
$ws.Range("G20").Value = 134.1
$ws.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=91014"
$ws.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202408/huvUNWz51724142842741.jpeg"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "2024-09-22"
$ws.Range("B21").NumberFormat = "General"
$ws.Range("C21").Value = "北京·《喜剧奇妙夜》一年一度喜剧大赛编剧团队编创/切西娅剧组演绎"
$ws.Range("D21").Value = "复兴路69号号3号楼6层601 爱乐汇艺术空间(五棵松万达广场)"
$ws.Range("E21").Value = "2024.09.22 19:30-09.22 21:00"
$ws.Range("F21").Value = 6
$ws.Range("G21").Value = 80
$ws.Range("H21").Value = "https://show.bilibili.com/platform/detail.html?id=90700"
$ws.Range("I21").Value = "//i1.hdslb.com/bfs/openplatform/202408/tCGETFGW1723613408321.jpeg"
$ws.Range("C22").Value = "北京·地狱双ip同人ONLY展"
$ws.Range("D22").Value = "双桥中路50号院 E50艺术园区"
$ws.Range("E22").Value = "2024.09.22 10:30-09.22 16:00"
$ws.Range("F22").Value = 134
$ws.Range("G22").Value = 105
$ws.Range("H22").Value = "https://show.bilibili.com/platform/detail.html?id=90931"
$ws.Range("I22").Value = "//i1.hdslb.com/bfs/openplatform/202408/c6ObwO4C1724055713128.jpeg"
$ws.Range("C23").Value = "北京·次元音浪Million Live⏤番音集结"
$ws.Range("D23").Value = "学清路38号金码大厦B座 北京想象空间"
$ws.Range("E23").Value = "2024.09.22 13:00-09.22 16:00"
$ws.Range("F23").Value = 46
$ws.Range("G23").Value = 88
$ws.Range("H23").Value = "https://show.bilibili.com/platform/detail.html?id=90657"
$ws.Range("I23").Value = "//i1.hdslb.com/bfs/openplatform/202408/Fn9CSOmf1723477511986.jpeg"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "2024-10-01"
$ws.Range("B24").NumberFormat = "General"
$ws.Range("C24").Value = "北京·IDO动漫游戏嘉年华47th"
$ws.Range("D24").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws.Range("E24").Value = "2024.10.01 09:30-10.02 17:00"
$ws.Range("F24").Value = 6910
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = "https://show.bilibili.com/platform/detail.html?id=83826"
$ws.Range("I24").Value = "//i0.hdslb.com/bfs/openplatform/202405/JL6boAFV1716882961702.jpeg"
$ws.Range("C25").Value = "北京·第19届IJOY漫展xCGF游戏节"
$ws.Range("D25").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E25").Value = "2024.10.01 09:00-10.02 17:00"
$ws.Range("F25").Value = 7472
$ws.Range("G25").Value = 8.800000000000001
$ws.Range("H25").Value = "https://show.bilibili.com/platform/detail.html?id=84127"
$ws.Range("I25").Value = "//i0.hdslb.com/bfs/openplatform/202405/iR6rV5311717039317028.jpeg"
$ws.Range("C26").Value = "北京·第19届IJOY漫展【Pile专场见面会】"
$ws.Range("E26").Value = "2024.10.01 14:50-10.01 16:30"
$ws.Range("F26").Value = 30
$ws.Range("G26").Value = 458
$ws.Range("H26").Value = "https://show.bilibili.com/platform/detail.html?id=91560"
$ws.Range("I26").Value = "//i1.hdslb.com/bfs/openplatform/202408/mBtVCKBp1724927832154.jpeg"
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "2024-10-02"
$ws.Range("B27").NumberFormat = "General"
$ws.Range("C27").Value = "北京·人气声优 青山渚 专场活动"
$ws.Range("D27").Value = "亦庄荣昌东街6号 北京亦创国际会展中心"
$ws.Range("E27").Value = "2024.10.02 11:50-10.02 15:40"
$ws.Range("F27").Value = 234
$ws.Range("H27").Value = "https://show.bilibili.com/platform/detail.html?id=91249"
$ws.Range("I27").Value = "//i2.hdslb.com/bfs/openplatform/202408/xHqpdFa41724641733192.png"
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "2024-10-02"
$ws.Range("B28").NumberFormat = "General"
$ws.Range("C28").Value = "北京·第19届IJOY漫展【217专场见面会】"
$ws.Range("D28").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E28").Value = "2024.10.02 12:25-10.02 16:30"
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 168
$ws.Range("H28").Value = "https://show.bilibili.com/platform/detail.html?id=91561"
$ws.Range("I28").Value = "//i0.hdslb.com/bfs/openplatform/202408/cAghXlck1724928163645.jpeg"
$ws.Range("C29").Value = "北京·第19届IJOY漫展【银发娘专场见面会】"
$ws.Range("D29").Value = "天辰东路7号 北京国家会议中心"
$ws.Range("E29").Value = "2024.10.02 12:25-10.02 16:30"
$ws.Range("F29").Value = 5
$ws.Range("G29").Value = 168
$ws.Range("H29").Value = "https://show.bilibili.com/platform/detail.html?id=91563"
$ws.Range("I29").Value = "//i0.hdslb.com/bfs/openplatform/202408/9Td79pPE1724928479521.jpeg"
$ws.Range("F31").Value = 1360
$ws.Range("F34").Value = 273
$ws.Range("F37").Value = 664
$ws.Range("F43").Value = 302
$ws.Range("F45").Value = 178
$ws.Range("F46").Value = 76
$ws.Range("F47").Value = 105
$ws.Range("F49").Value = 273
$ws.Range("F50").Value = 11
